# Updates cryptocurrency Price (D) and Volume(1h) (E) columns to the latest
# scraped figures, matching a new run of the GitHub Actions scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Assigning a numeric-looking string via .Value lets Excel
    # auto-coerce it to a number; force text (like the inline
    # strings already on the sheet) by flipping to a text format
    # for the write, then restoring the default "Normal" style.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.834.15"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.640.98"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.33%  "
Set-TextValue "D5" "218.32"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -0.45%  "
Set-TextValue "D9" "0.0622"
$ws.Range("E9").Value = "  -1.23%  "
Set-TextValue "D10" "19.25"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.871.15"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.642.25"
$ws.Range("E13").Value = "  -1.00%  "
Set-TextValue "D14" "4.15"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("E15").Value = "  -0.04%  "
Set-TextValue "D16" "65.29"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "26.837.83"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  -1.11%  "
Set-TextValue "D19" "216.32"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("E20").Value = "  -0.23%  "
Set-TextValue "D21" "4.36"
$ws.Range("E21").Value = "  -0.29%  "
Set-TextValue "D22" "6.57"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("E23").Value = "  -3.25%  "
Set-TextValue "D24" "9.20"
$ws.Range("E24").Value = "  -1.66%  "
Set-TextValue "D25" "147.51"
$ws.Range("E25").Value = "  +1.64%  "
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +1.40%  "
Set-TextValue "D33" "2.99"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "1.283.35"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -0.88%  "
Set-TextValue "D38" "0.533"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -0.31%  "
Set-TextValue "D41" "0.803"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "1.782.20"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("E44").Value = "  -6.12%  "
Set-TextValue "D45" "92.64"
$ws.Range("E45").Value = "  +1.14%  "
Set-TextValue "D46" "61.16"
$ws.Range("E46").Value = "  -1.22%  "
Set-TextValue "D47" "1.60"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -1.64%  "
Set-TextValue "D49" "7.58"
$ws.Range("E49").Value = "  -1.55%  "
Set-TextValue "D50" "0.0967"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  -0.26%  "
